# Quarterly financials update: add two new quarter columns.
#
# The "Period Ending" header rows (7, 38, 80) get two brand-new quarter
# columns inserted at D:E (2018-12-31 / 43465 and 2018-09-30 / 43373),
# with the previously-existing D:K values shifting right to F:M.
#
# Every other data row (8-35, 39-77, 81-102) just grows by two columns at
# the end (L, M), duplicating the rightmost existing column (K)'s value
# and formatting - matching how the source workbook was actually edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlShiftToRight = -4161

$headerRows = @(7, 38, 80)
$dataRowRanges = @(, (8, 35)), (, (39, 77)), (, (81, 102))

foreach ($row in $headerRows) {
    # Grab the existing (already correctly-formatted) date cell before we
    # start overwriting anything in this row.
    $fmtSource = $ws.Cells.Item($row, 4)

    # Shift existing D:K (cols 4-11) out to F:M (cols 6-13), working from
    # the rightmost column back so we never overwrite a cell before it has
    # been copied away.
    for ($c = 11; $c -ge 4; $c--) {
        $src = $ws.Cells.Item($row, $c)
        $dst = $ws.Cells.Item($row, $c + 2)
        $src.Copy($dst) | Out-Null
    }

    # Populate the two new leading quarter columns, reusing the original
    # date cell's formatting (copy-then-overwrite keeps the existing style
    # instead of minting a new one).
    $d1 = $ws.Cells.Item($row, 4)
    $d2 = $ws.Cells.Item($row, 5)
    $fmtSource.Copy($d1) | Out-Null
    $fmtSource.Copy($d2) | Out-Null
    $d1.Value2 = 43465
    $d2.Value2 = 43373
}

foreach ($bounds in $dataRowRanges) {
    $startRow = $bounds[0][0]
    $endRow = $bounds[0][1]
    for ($row = $startRow; $row -le $endRow; $row++) {
        $lastCol = $ws.Cells.Item($row, 11)
        $lastCol.Copy($ws.Cells.Item($row, 12)) | Out-Null
        $lastCol.Copy($ws.Cells.Item($row, 13)) | Out-Null
    }
}
